# Update "Sheet1" (the "bets to hit" odds-tracking sheet) with the current
# (Week 11) odds rows, appended after the existing Week 10 data (rows 2-134).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Week 11 rows: week, game, total_line, spread_line
$rows = @(
    @(135, 11, "NYJ_NE",  42.5,  5.5),
    @(136, 11, "WAS_MIA", 48.5, -3),
    @(137, 11, "GB_NYG",  43.5, -4.5),
    @(138, 11, "TB_BUF",  49.5,  5.5),
    @(139, 11, "CIN_PIT", 48.5, -1.5),
    @(140, 11, "HOU_TEN", 43.5, -2.5),
    @(141, 11, "CHI_MIN", 45.5,  2.5),
    @(142, 11, "CAR_ATL", 45.5,  3),
    @(143, 11, "LAC_JAX", 46.5, -1.5),
    @(144, 11, "SEA_LA",  45.5,  4.5),
    @(145, 11, "SF_ARI",  47.5, -1.5),
    @(146, 11, "BAL_CLE", 44.5, -8.5),
    @(147, 11, "KC_DEN",  45.5, -1.5),
    @(148, 11, "DET_PHI", 47.5,  3.5),
    @(149, 11, "DAL_LV",  45.5, -1.5)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Move the view/selection down to the newly-added rows.
$ws.Range("E138").Select() | Out-Null
